# Updated cryptos list on Mon Aug 28 23:18:48 UTC 2023 with GitHub Actions
#
# This script applies the per-cell market-data refresh produced by the
# scheduled scraper run: updated Price (column D) and Volume(1h) (column E)
# figures for each coin row, plus a swap in the ranking order between
# Polkadot and WrappedEther (rows 12 and 13 traded places, each carrying its
# own Coin name / Link / Price / Volume(1h) along for the ride).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") hold plain text in the
# source sheet (values like "26.162.82" or "0.5270" are display strings,
# not numbers). Excel's COM layer auto-coerces single-dot numeric-looking
# text typed into a Range.Value into a real number (and mints a new cell
# style along the way), so force the destination cell to Text format before
# writing, then restore its original (unstyled) look by copying the style
# of the untouched "Coin" cell in the same row.
function Set-TextValue {
    param($row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $plainStyle = $ws.Cells.Item($row, 2).Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $plainStyle
}

Set-TextValue 2 4 '26.162.82'
Set-TextValue 2 5 '  -0.49%  '
Set-TextValue 3 4 '1.657.39'
Set-TextValue 3 5 '  -0.73%  '
Set-TextValue 4 4 '1.004'
Set-TextValue 4 5 '  -0.36%  '
Set-TextValue 5 5 '  -0.26%  '
Set-TextValue 6 4 '0.5270'
Set-TextValue 6 5 '  -0.39%  '
Set-TextValue 7 5 '  -0.37%  '
Set-TextValue 8 4 '0.2685'
Set-TextValue 8 5 '  +1.07%  '
Set-TextValue 9 4 '0.06385'
Set-TextValue 9 5 '  +0.14%  '
Set-TextValue 10 4 '20.63'
Set-TextValue 10 5 '  -1.65%  '
Set-TextValue 11 4 '0.07699'
Set-TextValue 11 5 '  -1.80%  '
Set-TextValue 12 2 'WrappedEther'
Set-TextValue 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 12 4 '1.793.97'
Set-TextValue 12 5 '  +7.61%  '
Set-TextValue 13 2 'Polkadot'
Set-TextValue 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 13 4 '4.614'
Set-TextValue 13 5 '  +1.78%  '
Set-TextValue 14 4 '1.886.04'
Set-TextValue 14 5 '  -0.65%  '
Set-TextValue 15 4 '0.5645'
Set-TextValue 15 5 '  +0.64%  '
Set-TextValue 16 4 '0.0₅8270'
Set-TextValue 16 5 '  +1.99%  '
Set-TextValue 17 4 '65.77'
Set-TextValue 17 5 '  -0.13%  '
Set-TextValue 18 4 '26.148.95'
Set-TextValue 18 5 '  -0.43%  '
Set-TextValue 19 5 '  -0.38%  '
Set-TextValue 20 4 '4.688'
Set-TextValue 20 5 '  -0.92%  '
Set-TextValue 21 4 '10.37'
Set-TextValue 21 5 '  +0.83%  '
Set-TextValue 22 4 '191.00'
Set-TextValue 22 5 '  -4.91%  '
Set-TextValue 23 5 '  -1.29%  '
Set-TextValue 24 5 '  -0.31%  '
Set-TextValue 25 4 '147.07'
Set-TextValue 25 5 '  +0.62%  '
Set-TextValue 26 4 '0.1205'
Set-TextValue 26 5 '  -1.23%  '
Set-TextValue 27 4 '7.289'
Set-TextValue 27 5 '  +0.59%  '
Set-TextValue 28 4 '16.08'
Set-TextValue 28 5 '  -1.12%  '
Set-TextValue 29 4 '1.528'
Set-TextValue 29 5 '  -0.06%  '
Set-TextValue 30 4 '0.05644'
Set-TextValue 30 5 '  -4.55%  '
Set-TextValue 31 5 '  -0.36%  '
Set-TextValue 32 4 '3.497'
Set-TextValue 32 5 '  -0.64%  '
Set-TextValue 33 4 '3.382'
Set-TextValue 33 5 '  +1.37%  '
Set-TextValue 34 4 '1.582'
Set-TextValue 34 5 '  -1.07%  '
Set-TextValue 35 4 '2.798'
Set-TextValue 35 5 '  -0.81%  '
Set-TextValue 36 4 '0.9511'
Set-TextValue 36 5 '  -1.41%  '
Set-TextValue 37 5 '  -0.92%  '
Set-TextValue 38 4 '0.5781'
Set-TextValue 38 5 '  -0.54%  '
Set-TextValue 39 4 '0.01605'
Set-TextValue 39 5 '  -0.53%  '
Set-TextValue 40 4 '5.987'
Set-TextValue 40 5 '  +0.00%  '
Set-TextValue 41 5 '  -0.35%  '
Set-TextValue 42 4 '0.8351'
Set-TextValue 42 5 '  -2.62%  '
Set-TextValue 43 4 '1.029.07'
Set-TextValue 43 5 '  -4.54%  '
Set-TextValue 44 4 '101.60'
Set-TextValue 44 5 '  -1.28%  '
Set-TextValue 45 4 '1.795.87'
Set-TextValue 46 4 '58.56'
Set-TextValue 46 5 '  -0.03%  '
Set-TextValue 47 5 '  +3.54%  '
Set-TextValue 48 4 '0.05360'
Set-TextValue 48 5 '  +4.21%  '
Set-TextValue 49 4 '1.005'
Set-TextValue 49 5 '  -0.86%  '
Set-TextValue 50 4 '8.038'
Set-TextValue 50 5 '  -0.68%  '
Set-TextValue 51 4 '0.4343'
Set-TextValue 51 5 '  -1.54%  '
